$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.010.19"
$ws.Range("E2").Value = "  +0.28%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.124.77"
$ws.Range("E3").Value = "  +1.22%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'577.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'173.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.85%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.522"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.43%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  -2.48%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.154"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.49%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -0.01%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  -0.78%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "'37.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.09%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -1.23%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.641.00"
$ws.Range("E15").Value = "  +1.12%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "66.966.76"
$ws.Range("E16").Value = "  +0.16%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  -0.18%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.123.60"
$ws.Range("E18").Value = "  +1.08%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'16.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.27%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'476.60"
$ws.Range("D20").Style = "Normal"

# Row 21 - Polygon
$ws.Range("E21").Value = "  -0.36%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'7.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.55%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'84.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.71%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "'13.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.42%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  -2.42%  "

# Row 26 - RenderToken
$ws.Range("D26").Value = "'10.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.31%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.01%  "

# Row 28 - NEARProtocol
$ws.Range("D28").Value = "'7.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.98%  "

# Row 29 - ImmutableX
$ws.Range("E29").Value = "  -0.95%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value = "'28.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.38%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  +0.53%  "

# Row 33 - PEPE
$ws.Range("D33").Value = "0.0₃0954"
$ws.Range("E33").Value = "  -7.18%  "

# Row 34 - FirstDigitalUSD
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.01%  "

# Row 35 - Filecoin
$ws.Range("E35").Value = "  -0.43%  "

# Row 36 - Mantle
$ws.Range("D36").Value = "'0.977"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.69%  "

# Row 37 - Arweave
$ws.Range("D37").Value = "'47.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.39%  "

# Row 38 - row38(OKB->Stacks)
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'2.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.05%  "

# Row 39 - row39(Stacks->OKB)
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'50.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "

# Row 40 - TheGraph
$ws.Range("E40").Value = "  -1.42%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  +1.28%  "

# Row 42 - Cosmos
$ws.Range("E42").Value = "  -0.03%  "

# Row 43 - Maker
$ws.Range("D43").Value = "2.814.06"
$ws.Range("E43").Value = "  +1.51%  "

# Row 44 - Bittensor
$ws.Range("D44").Value = "'383.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.17%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  -1.61%  "

# Row 46 - dogwifhat
$ws.Range("D46").Value = "'2.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.42%  "

# Row 47 - Monero
$ws.Range("D47").Value = "'135.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.24%  "

# Row 48 - USDe
$ws.Range("E48").Value = "  -0.01%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").Value = "'24.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.55%  "

# Row 50 - ThetaToken
$ws.Range("E50").Value = "  -1.37%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -0.42%  "
